$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 3).Value = 188 + $r
    $ws.Rows.Item($r).RowHeight = 13.8
}

$ws.Range("C17").Select()
